# Updates the cryptos list (Coin / Link / Price / Volume(1h)) to the values
# published by the "Updated cryptos list" GitHub Actions run.
# Only columns B (Coin), C (Link), D (Price) and E (Volume(1h)) ever change;
# column A (the row index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, Coin, Link, Price, Volume(1h)
$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '46.088.70', '  -0.84%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.601.43', '  -0.37%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.05%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '311.28', '  +1.10%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '98.75', '  -2.41%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.599', '  -0.99%  '),
    @(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.01%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.582', '  +0.48%  '),
    @(10, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '39.01', '  +0.33%  '),
    @(11, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '54.51', '  -1.58%  '),
    @(12, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.0841', '  -0.22%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '8.14', '  -1.19%  '),
    @(14, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.001.67', '  -0.31%  '),
    @(15, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.107', '  +1.13%  '),
    @(16, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.607.48', '  -0.10%  '),
    @(17, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.918', '  +1.32%  '),
    @(18, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '14.88', '  -0.12%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '46.216.80', '  -0.86%  '),
    @(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000102', '  +0.68%  '),
    @(21, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '12.82', '  -3.86%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.72', '  -0.03%  '),
    @(23, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '295.02', '  +14.34%  '),
    @(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '72.61', '  +1.80%  '),
    @(25, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.07', '  +1.36%  '),
    @(26, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.28', '  +1.85%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '29.87', '  +5.55%  '),
    @(28, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.14%  '),
    @(29, 'LEO', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', '4.05', '  +0.97%  '),
    @(30, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.78', '  +2.48%  '),
    @(31, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '38.47', '  -3.66%  '),
    @(32, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.21', '  -2.68%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.26', '  +0.63%  '),
    @(34, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '3.58', '  -4.29%  '),
    @(35, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '155.60', '  +3.12%  '),
    @(36, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0840', '  +0.50%  '),
    @(37, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '2.21', '  -5.11%  '),
    @(38, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.79', '  -5.59%  '),
    @(39, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.121', '  +3.42%  '),
    @(40, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.123', '  +0.88%  '),
    @(41, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '22.34', '  +18.62%  '),
    @(42, 'Celestia', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia', '15.83', '  +0.04%  '),
    @(43, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0332', '  +2.27%  '),
    @(44, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '3.60', '  -0.85%  '),
    @(45, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '3.95', '  -5.58%  '),
    @(46, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.098.27', '  +2.04%  '),
    @(47, 'BitcoinSV', 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv', '98.17', '  +7.44%  '),
    @(48, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.998', '  -0.06%  '),
    @(49, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '9.62', '  +4.06%  '),
    @(50, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.202', '  +0.18%  '),
    @(51, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '108.50', '  -1.46%  ')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $coin = $r[1]
    $link = $r[2]
    $price = $r[3]
    $volume = $r[4]

    $ws.Range("B$rowNum").Value = $coin
    $ws.Range("C$rowNum").Value = $link

    # The Price column holds plain text (e.g. "46.088.70", "1.00", "0.0840")
    # that must never be reinterpreted as a number (which would collapse
    # trailing zeros or mis-parse multi-dot "thousands" separators).
    # Force a text number format before assigning, then restore the
    # cell's style so no stray formatting is introduced.
    $priceCell = $ws.Range("D$rowNum")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"

    $ws.Range("E$rowNum").Value = $volume
}
